$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2573958041829485
$ws.Range("C2").Value = 0.04748089165424574
$ws.Range("D2").Value = 0.03165052555379333
$ws.Range("E2").Value = 0.1647267180267136
$ws.Range("F2").Value = 0.7170521533020562
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("I2").Value = 0.5570539027569694
$ws.Range("K2").Value = 0.2750332023945532
$ws.Range("M2").Value = 0.2202528939603496
$ws.Range("O2").Value = 2.47108682482488

$ws.Range("B3").Value = 0.2256133626913197
$ws.Range("C3").Value = 0.04230195993052632
$ws.Range("D3").Value = 0.02950460241034136
$ws.Range("E3").Value = 0.153940755802644
$ws.Range("F3").Value = 0.7163095727574884
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("I3").Value = 0.5619326161867058
$ws.Range("K3").Value = 0.2399888997760513
$ws.Range("M3").Value = 0.1977350549761994
$ws.Range("O3").Value = 2.481820132211851

$ws.Range("B4").Value = 0.2060647825611284
$ws.Range("C4").Value = 0.03910215581566945
$ws.Range("D4").Value = 0.02817525147609246
$ws.Range("E4").Value = 0.1474358730610774
$ws.Range("F4").Value = 0.7163156781840172
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("I4").Value = 0.5652477437201249
$ws.Range("K4").Value = 0.2184091649824609
$ws.Range("M4").Value = 0.1839742000649167
$ws.Range("O4").Value = 2.48992116368197

$ws.Range("B5").Value = 0.198090472907694
$ws.Range("C5").Value = 0.0377932470180582
$ws.Range("D5").Value = 0.02763060759778568
$ws.Range("E5").Value = 0.1448144618963596
$ws.Range("F5").Value = 0.7164343057102229
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("I5").Value = 0.5666789925479563
$ws.Range("K5").Value = 0.2096000144404826
$ws.Range("M5").Value = 0.1783829147458533
$ws.Range("O5").Value = 2.49360203497055

$ws.Range("B6").Value = 0.1967658687564438
$ws.Range("C6").Value = 0.03757560534734239
$ws.Range("D6").Value = 0.02753999429189236
$ws.Range("E6").Value = 0.1443809465394281
$ws.Range("F6").Value = 0.7164610168747743
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("I6").Value = 0.566921499146261
$ws.Range("K6").Value = 0.2081363550997537
$ws.Range("M6").Value = 0.1774554746049688
$ws.Range("O6").Value = 2.494236162507875

$ws.Range("B7").Value = 0.2059572704447419
$ws.Range("C7").Value = 0.03908452344474256
$ws.Range("D7").Value = 0.02816791801046037
$ws.Range("E7").Value = 0.1474004011273919
$ws.Range("F7").Value = 0.7163168078561455
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("I7").Value = 0.5652667209518931
$ws.Range("K7").Value = 0.2182904227715881
$ws.Range("M7").Value = 0.183898727756663
$ws.Range("O7").Value = 2.489969268410391

$ws.Range("B8").Value = 0.2464445435215907
$ws.Range("C8").Value = 0.04569935633928424
$ws.Range("D8").Value = 0.03091306692943618
$ws.Range("E8").Value = 0.1609831416518688
$ws.Range("F8").Value = 0.7167001816052903
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("I8").Value = 0.5586697228257549
$ws.Range("K8").Value = 0.2629631383053095
$ws.Range("M8").Value = 0.2124751952888957
$ws.Range("O8").Value = 2.47447396015481

$ws.Range("B9").Value = 0.3255538798997577
$ws.Range("C9").Value = 0.05851154655238133
$ws.Range("D9").Value = 0.03620197935609326
$ws.Range("E9").Value = 0.1885647269042394
$ws.Range("F9").Value = 0.7211214530595029
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("I9").Value = 0.5482712584070697
$ws.Range("K9").Value = 0.3500553480225221
$ws.Range("M9").Value = 0.2690345788772888
$ws.Range("O9").Value = 2.456088234992194

$ws.Range("B10").Value = 0.3834852254432519
$ws.Range("C10").Value = 0.06782640898629211
$ws.Range("D10").Value = 0.04002908673437844
$ws.Range("E10").Value = 0.2094237722106271
$ws.Range("F10").Value = 0.7266129998833577
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("I10").Value = 0.5421825645107319
$ws.Range("K10").Value = 0.4137155257719769
$ws.Range("M10").Value = 0.3109158923864754
$ws.Range("O10").Value = 2.449917869584255

$ws.Range("B11").Value = 0.409795378727523
$ws.Range("C11").Value = 0.07204248513725986
$ws.Range("D11").Value = 0.04175717665963674
$ws.Range("E11").Value = 0.2190463014861734
$ws.Range("F11").Value = 0.7295997812108084
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("I11").Value = 0.5397502699875112
$ws.Range("K11").Value = 0.4426025116473795
$ws.Range("M11").Value = 0.3300421186969231
$ws.Range("O11").Value = 2.448709000461747

$ws.Range("B12").Value = 0.4197517723568183
$ws.Range("C12").Value = 0.07363590661034891
$ws.Range("D12").Value = 0.04240968049591487
$ws.Range("E12").Value = 0.2227096001105764
$ws.Range("F12").Value = 0.7308011475503662
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("I12").Value = 0.5388778192227868
$ws.Range("K12").Value = 0.4535304850825526
$ws.Range("M12").Value = 0.3372954822992753
$ws.Range("O12").Value = 2.448481378540009

$ws.Range("B13").Value = 0.4176077902129407
$ws.Range("C13").Value = 0.07329287415457486
$ws.Range("D13").Value = 0.0422692365560593
$ws.Range("E13").Value = 0.2219197741283665
$ws.Range("F13").Value = 0.7305392824520851
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("I13").Value = 0.5390635541028352
$ws.Range("K13").Value = 0.4511774432106677
$ws.Range("M13").Value = 0.3357328657208072
$ws.Range("O13").Value = 2.448520159406854

$ws.Range("B14").Value = 0.410614634119213
$ws.Range("C14").Value = 0.07217363960768353
$ws.Range("D14").Value = 0.04181089654312586
$ws.Range("E14").Value = 0.219347291987539
$ws.Range("F14").Value = 0.7296972082864315
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("I14").Value = 0.5396775183547646
$ws.Range("K14").Value = 0.4435017839176396
$ws.Range("M14").Value = 0.330638643273069
$ws.Range("O14").Value = 2.448685659302669

$ws.Range("B15").Value = 0.4063302382385814
$ws.Range("C15").Value = 0.07148766892133551
$ws.Range("D15").Value = 0.04152990338573659
$ws.Range("E15").Value = 0.2177741131480744
$ws.Range("F15").Value = 0.7291905755645089
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("I15").Value = 0.5400599214475577
$ws.Range("K15").Value = 0.4387987865819127
$ws.Range("M15").Value = 0.3275196755406284
$ws.Range("O15").Value = 2.448817015564657

$ws.Range("B16").Value = 0.3817648855556115
$ws.Range("C16").Value = 0.06755044540766164
$ws.Range("D16").Value = 0.0399158898858758
$ws.Range("E16").Value = 0.2087976268839071
$ws.Range("F16").Value = 0.7264276455719383
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("I16").Value = 0.5423483182619222
$ws.Range("K16").Value = 0.4118261948931377
$ws.Range("M16").Value = 0.3096674422559715
$ws.Range("O16").Value = 2.450029059170333

$ws.Range("B17").Value = 0.3666834462312352
$ws.Range("C17").Value = 0.06512959344172486
$ws.Range("D17").Value = 0.0389224215774675
$ws.Range("E17").Value = 0.2033252420597549
$ws.Range("F17").Value = 0.7248578823874468
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("I17").Value = 0.5438386677224152
$ws.Range("K17").Value = 0.3952604986911581
$ws.Range("M17").Value = 0.2987346857489257
$ws.Range("O17").Value = 2.451182160968472

$ws.Range("B18").Value = 0.3580049731122585
$ws.Range("C18").Value = 0.0637351815080649
$ws.Range("D18").Value = 0.03834979454769183
$ws.Range("E18").Value = 0.2001902361395054
$ws.Range("F18").Value = 0.7240009859789751
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("I18").Value = 0.5447276382183475
$ws.Range("K18").Value = 0.3857255682301002
$ws.Range("M18").Value = 0.2924534377089572
$ws.Range("O18").Value = 2.451995790778568

$ws.Range("B19").Value = 0.3550659144171107
$ws.Range("C19").Value = 0.06326271534368288
$ws.Range("D19").Value = 0.0381557060405342
$ws.Range("E19").Value = 0.1991309278267295
$ws.Range("F19").Value = 0.723718752233772
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("I19").Value = 0.5450340809793914
$ws.Range("K19").Value = 0.3824960553068024
$ws.Range("M19").Value = 0.2903279144330071
$ws.Range("O19").Value = 2.452297091811971

$ws.Range("B20").Value = 0.3682893116560138
$ws.Range("C20").Value = 0.0653875048328274
$ws.Range("D20").Value = 0.03902830348704356
$ws.Range("E20").Value = 0.2039064846067049
$ws.Range("F20").Value = 0.7250202260695957
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("I20").Value = 0.5436767297859362
$ws.Range("K20").Value = 0.3970246506635817
$ws.Range("M20").Value = 0.2998977734480235
$ws.Range("O20").Value = 2.451043843320434

$ws.Range("B21").Value = 0.4126688790654214
$ws.Range("C21").Value = 0.07250247087846162
$ws.Range("D21").Value = 0.04194557358052009
$ws.Range("E21").Value = 0.2201023624359593
$ws.Range("F21").Value = 0.7299426364013542
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("I21").Value = 0.5394958623971249
$ws.Range("K21").Value = 0.4457566109230697
$ws.Range("M21").Value = 0.332134650079567
$ws.Range("O21").Value = 2.448630798958959

$ws.Range("B22").Value = 0.4416342682096399
$ws.Range("C22").Value = 0.0771343316283577
$ws.Range("D22").Value = 0.04384116759167256
$ws.Range("E22").Value = 0.2308008174482126
$ws.Range("F22").Value = 0.7335697260520675
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("I22").Value = 0.5370467745601069
$ws.Range("K22").Value = 0.4775419941501582
$ws.Range("M22").Value = 0.3532655477519171
$ws.Range("O22").Value = 2.448395390555277

$ws.Range("B23").Value = 0.4261786458274059
$ws.Range("C23").Value = 0.07466390100798037
$ws.Range("D23").Value = 0.04283047211082902
$ws.Range("E23").Value = 0.2250803857136816
$ws.Range("F23").Value = 0.7315963412949174
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("I23").Value = 0.5383279469125561
$ws.Range("K23").Value = 0.4605835357412786
$ws.Range("M23").Value = 0.341981893494868
$ws.Range("O23").Value = 2.448398155227068

$ws.Range("B24").Value = 0.3675633246329255
$ws.Range("C24").Value = 0.06527091128747031
$ws.Range("D24").Value = 0.03898043884797886
$ws.Range("E24").Value = 0.2036436701688302
$ws.Range("F24").Value = 0.724946688512297
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("I24").Value = 0.5437498417900564
$ws.Range("K24").Value = 0.3962271120681748
$ws.Range("M24").Value = 0.29937192860978
$ws.Range("O24").Value = 2.451105907339979

$ws.Range("B25").Value = 0.3041849854042198
$ws.Range("C25").Value = 0.05506269063980085
$ws.Range("D25").Value = 0.03478140592350343
$ws.Range("E25").Value = 0.1810000314860218
$ws.Range("F25").Value = 0.7195318920116591
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("I25").Value = 0.550812202364348
$ws.Range("K25").Value = 0.3265508511610449
$ws.Range("M25").Value = 0.2536770058462992
$ws.Range("O25").Value = 2.459775065628975
